$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cells F10, F11, F12, F13, F17 from "Necessary" to "Complete"
$ws.Range("F10").Value = "Complete"
$ws.Range("F11").Value = "Complete"
$ws.Range("F12").Value = "Complete"
$ws.Range("F13").Value = "Complete"
$ws.Range("F17").Value = "Complete"

# Match the green fill style used by other "Complete" cells (e.g. F9, F14)
$ws.Range("F10").Interior.Color = $ws.Range("F9").Interior.Color
$ws.Range("F11").Interior.Color = $ws.Range("F9").Interior.Color
$ws.Range("F12").Interior.Color = $ws.Range("F9").Interior.Color
$ws.Range("F13").Interior.Color = $ws.Range("F9").Interior.Color
$ws.Range("F17").Interior.Color = $ws.Range("F9").Interior.Color

# Update the frozen pane top-left cell and active selection to reflect scrolled view
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("A9:F41").Select()
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("D15").Select()
